# Updated symbol list on Tue Dec 27 19:19:45 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores numeric-looking values as TEXT (they were
# written as inlineStr in the source file), so each Price update is entered
# with a leading apostrophe (forces text entry) and the cell style is then
# reset to "Normal" so no stray quote-prefix / number-format style sticks to
# the cell (keeping it identical to the untouched cells around it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText {
    param($addr, $text)
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).Style = "Normal"
}

# Straight price refreshes (coin/link/volume unchanged)
Set-PriceText "D2"  "245.87"
Set-PriceText "D4"  "5.333"
Set-PriceText "D5"  "0.05828"
Set-PriceText "D6"  "6.477"
Set-PriceText "D7"  "3.352"
Set-PriceText "D8"  "0.8111"
Set-PriceText "D9"  "0.9190"
Set-PriceText "D11" "0.07345"
Set-PriceText "D12" "0.03083"
Set-PriceText "D13" "0.03074"
Set-PriceText "D14" "0.09363"
Set-PriceText "D15" "3.878"
Set-PriceText "D16" "0.001564"

# Rows 18-24: ranking reshuffled (each row's Coin/Link/Volume shift to the
# row above; new Price values accompany every row; "Bestin24h" badge moves
# from NitroEx's row to One's row, which now lands at the bottom, row 24)
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-PriceText "D18" "0.006065"
$ws.Range("E18").Value = "17TigerCashTCH"

$ws.Range("B19").Value = "BitKan"
$ws.Range("C19").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-PriceText "D19" "0.001248"
$ws.Range("E19").Value = "18BitKanKAN"

$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-PriceText "D20" "0.004689"
$ws.Range("E20").Value = "19HotbitTokenHTB"

$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-PriceText "D21" "0.00008814"
$ws.Range("E21").Value = "20NitroExNTX"

$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-PriceText "D22" "3.593"
$ws.Range("E22").Value = "21LEOLEO"

$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-PriceText "D23" "2.158"
$ws.Range("E23").Value = "22BTSETokenBTSE"

$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-PriceText "D24" "0.01080"
$ws.Range("E24").Value = "23OneONEBestin24h"

Set-PriceText "D25" "0.3231"
Set-PriceText "D40" "0.03834"

# Rows 41-43: ranking reshuffled similarly ("Worstin24h" badge moves from
# KickToken's row to BOLO's row, row 48)
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-PriceText "D41" "0.006350"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-PriceText "D42" "0.1064"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-PriceText "D43" "0.003205"
$ws.Range("E43").Value = "42CEJICEJI"

Set-PriceText "D44" "0.007783"
Set-PriceText "D45" "0.00005258"
Set-PriceText "D47" "0.6811"
Set-PriceText "D48" "0.001860"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
